$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 524.96738300000004
$ws.Range("E2").Value = 803.06557899999996

$ws.Range("E2").Select()
